$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new total row: sum of the line-item totals (F5:F7) into F9,
# matching the existing currency formatting used by the column above.
$ws.Range("F9").Formula = "=SUM(F5:F7)"
$ws.Range("F9").NumberFormat = $ws.Range("F7").NumberFormat

# Move the active selection to the newly-added total cell.
$ws.Range("F9").Select()
